# CHARM v3 documentation: remove the standalone "Probs" function
# documentation section (its own Heading1 + description paragraphs).
# The Probs() function usage elsewhere in the document (bullet list,
# inline Probs() mentions, etc.) is left untouched - only the dedicated
# "Probs" chapter (describing M_it, fire_it, x_i, intervention_it,
# deathrate_t, debug, v.p.it, v.x, v.rr, modifyRisk(), ... "returned as
# v.p.") is deleted, merging back into the blank spacer paragraphs that
# precede the "References" section.

$d = $word.ActiveDocument

# --- Locate the unique sentence that opens the "Probs" chapter body ---
$rAfterHeading = $d.Content.Duplicate
$foundAfterHeading = $rAfterHeading.Find.Execute( `
    "This function calculates transition probabilities for each individual at each cycle given individual risk factors and previous health status. Its arguments are:", `
    $false, $true, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $foundAfterHeading) {
    throw "Could not find the 'This function calculates transition probabilities...' paragraph"
}

# --- Within the text before that, find the last "Probs" occurrence:    ---
# --- that is the lone heading run "Probs" (after the page break run). ---
$beforeHeading = $d.Range(0, $rAfterHeading.Start)
$probe = $beforeHeading.Duplicate
$headingStart = -1
$headingEnd = -1
while ($probe.Find.Execute("Probs", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)) {
    $headingStart = $probe.Start
    $headingEnd = $probe.End
    $probe.Collapse(0)
    $probe.End = $beforeHeading.End
}

if ($headingStart -eq -1) {
    throw "Could not find the 'Probs' heading run"
}

# --- Locate the end of the chapter: "...is returned as v.p." ---
$rEndSentence = $d.Content.Duplicate
$foundEndSentence = $rEndSentence.Find.Execute( `
    "The adjusted vector of transition probabilities for the individual is returned as", `
    $false, $true, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $foundEndSentence) {
    throw "Could not find the 'The adjusted vector of transition probabilities...' paragraph"
}

# From the end of that match (right before " v.p.") walk forward through
# " v.p." and then through the paragraph mark that closes this paragraph,
# so the whole paragraph (including its pilcrow) is removed.
$tail = $d.Range($rEndSentence.End, $rEndSentence.End + 12)
$tailText = $tail.Text
$parts = $tailText.Split([char]13)
$sentenceTailLen = $parts[0].Length
$chapterEnd = $rEndSentence.End + $sentenceTailLen + 1

# --- Delete in reverse document order so earlier offsets are not ---
# --- invalidated by a later deletion shifting the text around.   ---

# 1) Delete everything from the chapter body through the paragraph
#    mark that ends "...is returned as v.p." - this removes every
#    paragraph of the "Probs" chapter, merging what's left of the
#    (now textless) heading paragraph directly into the blank
#    spacer paragraphs that lead into "References".
$d.Range($rAfterHeading.Start, $chapterEnd).Delete()

# 2) Delete the heading's own "Probs" text run (keeps the leading
#    space + page-break runs, and the paragraph itself, in place).
$d.Range($headingStart, $headingEnd).Delete()
